# Generate Report for Handback
# The localization entry for "9c330023-7dc2-49c1-967b-e5ba9eec275f.md" has
# been handed back (in sync with en-US) for both the zh-cn and de-de
# locales. Update the Overview sheet's status columns plus each locale
# sheet's Status and Latest Handback DateTime cells.

$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"

# --- Overview sheet: row for the 9c330023... file, zh-cn (B) and de-de (C) columns
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $status
$overview.Range("C3").Value = $status

# --- zh-cn sheet: Status (B3) and Latest Handback DateTime (G3)
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = $status
$zhcn.Range("G3").Value = "2016-02-26 06:34:08"

# --- de-de sheet: Status (B3) and Latest Handback DateTime (G3)
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = $status
$dede.Range("G3").Value = "2016-02-26 06:34:32"
